# Adds a new "2022-Q4" sheet (fund-holdings detail) right after "总计",
# pushing every existing quarter sheet down by one tab position, and
# inserts the corresponding 2022-Q4 summary row at the top of "总计".

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item(1)

# --- 1. Insert the new "2022-Q4" worksheet right after "总计" --------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Use the "2022-Q3" sheet (looked up by name, since positional indices shift
# once the new tab is inserted) as a formatting donor so the new sheet's
# header/index-column styles match the existing quarter sheets without
# fabricating new style entries.
$formatDonor = $wb.Worksheets.Item("2022-Q3")

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# index, fund code, fund name, fund size, stock position, position ratio, market value, rank
$rows = @(
    @(0, "501029", "华宝标普中国A股红利机会指数（LOF）A", "10.97", "94.25", "1.81", "0.1986", 4),
    @(1, "005125", "华宝标普中国A股红利机会指数C",          "3.29",  "94.25", "1.81", "0.0595", 4),
    @(2, "512040", "富国中证价值ETF",                        "3.39",  "99.29", "1.07", "0.0363", 6),
    @(3, "501307", "银河中证沪港深高股息指数（LOF）A",       "0.16",  "93.15", "1.46", "0.0023", 5),
    @(4, "501308", "银河中证沪港深高股息指数（LOF）C",       "0.01",  "93.15", "1.46", "0.0001", 5)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = "'" + $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = "'" + $row[3]
    $newSheet.Cells.Item($r, 5).Value = "'" + $row[4]
    $newSheet.Cells.Item($r, 6).Value = "'" + $row[5]
    $newSheet.Cells.Item($r, 7).Value = "'" + $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
}

# Match header row + index column styling to the sibling quarter sheets
# (format-only paste re-uses the existing style entries).
$formatDonor.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$formatDonor.Range("A2").Copy()
$newSheet.Range("A2:A6").PasteSpecial(-4122)

# --- 2. Insert the 2022-Q4 row at the top of the "总计" summary sheet ------
$summaryData = @(
    @(0, "2022-Q4", 5, 0.3),
    @(1, "2022-Q3", 13, 2.3),
    @(2, "2022-Q2", 15, 2.69),
    @(3, "2022-Q1", 10, 3.23),
    @(4, "2021-Q4", 9, 2.72),
    @(5, "2021-Q3", 7, 2.46),
    @(6, "2021-Q2", 9, 3.31),
    @(7, "2021-Q1", 6, 0.74),
    @(8, "2020-Q4", 6, 0.8100000000000001)
)

# Row 10 is brand new - clone row 9's index-column style onto it first so the
# new "A10" cell matches the others before we overwrite the value.
$summary.Cells.Item(9, 1).Copy($summary.Cells.Item(10, 1))

for ($i = 0; $i -lt $summaryData.Length; $i++) {
    $row = $i + 2
    $entry = $summaryData[$i]
    $summary.Cells.Item($row, 1).Value = $entry[0]
    $summary.Cells.Item($row, 2).Value = $entry[1]
    $summary.Cells.Item($row, 3).Value = $entry[2]
    $summary.Cells.Item($row, 4).Value = $entry[3]
}

# --- 3. Restore the active tab to the last sheet (2020-Q4), matching the ---
# original workbook's selection state.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
